# "Added camera controls and README.txt."
# Mark several CAMERA & VIEW / LIGHTING features as completed on Milestone III.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows that are fully marked complete for Milestone III (Student column = III,
# Milestone III Complete column = X) -> these feed the G column formula and
# ripple into the Milestone III point totals.
$completedRows = @(12, 15, 24, 29, 31)
foreach ($r in $completedRows) {
    $ws.Cells.Item($r, 5).Value = "III"   # column E - Student(I, II, or III)
    $ws.Cells.Item($r, 6).Value = "X"     # column F - Milestone III Complete(X)
}

# Rows that are labeled as Milestone III features (camera controls, etc.) but
# not yet checked off as complete -> only column E changes.
$labelOnlyRows = @(57, 63, 68, 72, 82)
foreach ($r in $labelOnlyRows) {
    $ws.Cells.Item($r, 5).Value = "III"   # column E - Student(I, II, or III)
}

# Update the window view to reflect the new scroll position / selection.
$ws.Range("F15").Select()
